$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price (column D) cells whose new value reads as a plain number; force text
# storage (matching the source data, which stores every Price/Volume cell as
# a literal string) by temporarily switching to a text number format, then
# restore the default "Normal" style so no stray formatting is introduced.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "323.99"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "101.62"
$ws.Range("D6").Style = "Normal"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "40.12"
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0923"
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "8.46"
$ws.Range("D12").Style = "Normal"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.00"
$ws.Range("D13").Style = "Normal"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "16.43"
$ws.Range("D15").Style = "Normal"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "8.10"
$ws.Range("D18").Style = "Normal"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "76.01"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.70"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "266.27"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.31"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "10.16"
$ws.Range("D25").Style = "Normal"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "22.93"
$ws.Range("D28").Style = "Normal"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "175.83"
$ws.Range("D30").Style = "Normal"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.09"
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0902"
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "35.29"
$ws.Range("D33").Style = "Normal"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.07"
$ws.Range("D34").Style = "Normal"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.132"
$ws.Range("D35").Style = "Normal"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.58"
$ws.Range("D36").Style = "Normal"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.93"
$ws.Range("D38").Style = "Normal"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.79"
$ws.Range("D40").Style = "Normal"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.239"
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "70.24"
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "120.85"
$ws.Range("D44").Style = "Normal"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "91.61"
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "11.94"
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "5.53"
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.25"
$ws.Range("D49").Style = "Normal"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.101"
$ws.Range("D51").Style = "Normal"

# Remaining Price/Volume cells: plain text assignment is safe as-is.
$ws.Range("D2").Value = "42.568.61"
$ws.Range("E2").Value = "  -2.05%  "
$ws.Range("D3").Value = "2.357.20"
$ws.Range("E3").Value = "  +0.14%  "
$ws.Range("E4").Value = "  -0.21%  "
$ws.Range("E5").Value = "  +3.53%  "
$ws.Range("E6").Value = "  -7.09%  "
$ws.Range("E7").Value = "  -0.68%  "
$ws.Range("E9").Value = "  -2.31%  "
$ws.Range("E10").Value = "  -7.49%  "
$ws.Range("E11").Value = "  -1.73%  "
$ws.Range("E12").Value = "  -4.89%  "
$ws.Range("E13").Value = "  -2.95%  "
$ws.Range("E14").Value = "  -0.04%  "
$ws.Range("E15").Value = "  +0.05%  "
$ws.Range("D16").Value = "2.716.67"
$ws.Range("E16").Value = "  +0.38%  "
$ws.Range("D17").Value = "2.360.00"
$ws.Range("E17").Value = "  -2.52%  "
$ws.Range("E18").Value = "  +11.92%  "
$ws.Range("D19").Value = "42.601.00"
$ws.Range("E19").Value = "  -1.92%  "
$ws.Range("E20").Value = "  -1.89%  "
$ws.Range("E21").Value = "  +1.26%  "
$ws.Range("E22").Value = "  +7.33%  "
$ws.Range("E23").Value = "  +3.01%  "
$ws.Range("E24").Value = "  -9.71%  "
$ws.Range("E25").Value = "  +10.07%  "
$ws.Range("E26").Value = "  +0.22%  "
$ws.Range("E27").Value = "  -5.26%  "
$ws.Range("E28").Value = "  +1.65%  "
$ws.Range("E29").Value = "  -2.15%  "
$ws.Range("E30").Value = "  +1.42%  "
$ws.Range("E31").Value = "  -3.12%  "
$ws.Range("E32").Value = "  -3.15%  "
$ws.Range("E33").Value = "  -10.38%  "
$ws.Range("E34").Value = "  +0.57%  "
$ws.Range("E35").Value = "  -0.17%  "
$ws.Range("E36").Value = "  -7.91%  "
$ws.Range("E37").Value = "  -4.75%  "
$ws.Range("E38").Value = "  +6.63%  "
$ws.Range("E39").Value = "  +2.10%  "
$ws.Range("E40").Value = "  -8.77%  "
$ws.Range("E41").Value = "  +1.69%  "
$ws.Range("E42").Value = "  +2.18%  "
$ws.Range("E43").Value = "  -3.01%  "
$ws.Range("E44").Value = "  +8.93%  "
$ws.Range("E45").Value = "  -0.25%  "
$ws.Range("E46").Value = "  +22.34%  "
$ws.Range("E47").Value = "  -6.77%  "
$ws.Range("E48").Value = "  -1.92%  "
$ws.Range("E49").Value = "  -1.36%  "
$ws.Range("E50").Value = "  -3.84%  "
$ws.Range("E51").Value = "  +0.36%  "
